# "Correcta Autoevaluacio Ivan Rubi"
# Fixes the self-assessment sheet:
#  - marks the "Path Coverage" (column K) topic as worked on in the
#    "alumne" table (row 19), which recalculates the dependent totals
#  - lowers the max grade aimed for (C37) from 8 to 7, which also
#    recalculates the final grade (C39)
#  - adds an (empty, underlined) cell at H32, just under the second
#    table, matching a small layout tweak made while reviewing the sheet
#  - leaves the active selection on H32, where the author ended up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Path Coverage" (column K) as worked on for "Total alumne" row.
$ws.Range("K19").Value = 1

# Lower the targeted maximum grade for the assignment.
$ws.Range("C37").Value = 7

# Add the small underlined marker cell below the second table (creates
# row 32 with an empty, underlined H32 cell and its own style/font).
$h32 = $ws.Range("H32")
$h32.Value = $null
$h32.Font.Underline = 1

# Leave the selection where the author left it.
[void]$h32.Select()
